# Weekly refresh of fruit/vegetable price data
# (Fruta / hortaliza, semanal)
#
# The underlying source data for rows 3-14 got reshuffled between a
# weekly pull; this re-applies the new Fecha / Calidad / Volumen /
# Precio / Unidad de comercialización values per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (sourced from original row 9)
$ws.Range("D3").Value = 45030
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 16000
$ws.Range("P3").Value = 15500
$ws.Range("Q3").Value = '$/caja 18 kilos granel'
$ws.Range("S3").Value = 861

# Row 4 (sourced from original row 8)
$ws.Range("D4").Value = 44819
$ws.Range("N4").Value = 25000
$ws.Range("O4").Value = 26000
$ws.Range("P4").Value = 25500
$ws.Range("Q4").Value = '$/caja 18 kilos granel'
$ws.Range("S4").Value = 1417

# Row 5 (sourced from original row 13)
$ws.Range("D5").Value = 45014
$ws.Range("L5").Value = 'Primera'
$ws.Range("N5").Value = 13000
$ws.Range("O5").Value = 14000
$ws.Range("P5").Value = 13600
$ws.Range("S5").Value = 756

# Row 6 (sourced from original row 14)
$ws.Range("D6").Value = 45014
$ws.Range("L6").Value = 'Segunda'
$ws.Range("M6").Value = 20
$ws.Range("N6").Value = 10000
$ws.Range("O6").Value = 10000
$ws.Range("P6").Value = 10000
$ws.Range("S6").Value = 556

# Row 8 (sourced from original row 4)
$ws.Range("D8").Value = 44699
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 22000
$ws.Range("P8").Value = 21000
$ws.Range("Q8").Value = '$/caja 18 kilos'
$ws.Range("S8").Value = 1167

# Row 9 (sourced from original row 5)
$ws.Range("D9").Value = 44699
$ws.Range("L9").Value = 'Segunda'
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = 18000
$ws.Range("O9").Value = 18000
$ws.Range("P9").Value = 18000
$ws.Range("Q9").Value = '$/caja 18 kilos'
$ws.Range("S9").Value = 1000

# Row 10 (sourced from original row 6)
$ws.Range("D10").Value = 45002
$ws.Range("N10").Value = 12000
$ws.Range("O10").Value = 13000
$ws.Range("P10").Value = 12500
$ws.Range("S10").Value = 694

# Row 11 (sourced from original row 10)
$ws.Range("D11").Value = 44516
$ws.Range("N11").Value = 33000
$ws.Range("O11").Value = 34000
$ws.Range("P11").Value = 33500
$ws.Range("S11").Value = 1861

# Row 12 (sourced from original row 11)
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 100
$ws.Range("N12").Value = 14000
$ws.Range("O12").Value = 15000
$ws.Range("P12").Value = 14500
$ws.Range("S12").Value = 806

# Row 13 (sourced from original row 12)
$ws.Range("D13").Value = 44280
$ws.Range("L13").Value = 'Segunda'
$ws.Range("N13").Value = 12000
$ws.Range("O13").Value = 12000
$ws.Range("P13").Value = 12000
$ws.Range("S13").Value = 667

# Row 14 (sourced from original row 3)
$ws.Range("D14").Value = 44687
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = 18000
$ws.Range("O14").Value = 19000
$ws.Range("P14").Value = 18500
$ws.Range("S14").Value = 1028
